# Apply updated cryptocurrency price/volume data to the active worksheet.
# Generated to reproduce the commit "Updated cryptos list ... with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Col = 'D'; Row = 2; Value = '27.508.48' },
    @{ Col = 'E'; Row = 2; Value = '  -2.17%  ' },
    @{ Col = 'D'; Row = 3; Value = '1.844.20' },
    @{ Col = 'E'; Row = 3; Value = '  -2.79%  ' },
    @{ Col = 'D'; Row = 4; Value = '1.005' },
    @{ Col = 'E'; Row = 4; Value = '  -1.06%  ' },
    @{ Col = 'D'; Row = 5; Value = '333.81' },
    @{ Col = 'E'; Row = 5; Value = '  -1.13%  ' },
    @{ Col = 'D'; Row = 6; Value = '1.005' },
    @{ Col = 'E'; Row = 6; Value = '  -1.00%  ' },
    @{ Col = 'D'; Row = 7; Value = '0.4610' },
    @{ Col = 'E'; Row = 7; Value = '  -3.48%  ' },
    @{ Col = 'D'; Row = 8; Value = '0.3835' },
    @{ Col = 'E'; Row = 8; Value = '  -3.63%  ' },
    @{ Col = 'D'; Row = 9; Value = '46.50' },
    @{ Col = 'E'; Row = 9; Value = '  -2.42%  ' },
    @{ Col = 'D'; Row = 10; Value = '0.07892' },
    @{ Col = 'E'; Row = 10; Value = '  -2.22%  ' },
    @{ Col = 'D'; Row = 11; Value = '0.9812' },
    @{ Col = 'E'; Row = 11; Value = '  -4.45%  ' },
    @{ Col = 'D'; Row = 12; Value = '21.23' },
    @{ Col = 'E'; Row = 12; Value = '  -4.02%  ' },
    @{ Col = 'B'; Row = 13; Value = 'Polkadot' },
    @{ Col = 'C'; Row = 13; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot' },
    @{ Col = 'D'; Row = 13; Value = '5.909' },
    @{ Col = 'E'; Row = 13; Value = '  -2.20%  ' },
    @{ Col = 'B'; Row = 14; Value = 'WrappedEther' },
    @{ Col = 'C'; Row = 14; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' },
    @{ Col = 'D'; Row = 14; Value = '1.825.41' },
    @{ Col = 'E'; Row = 14; Value = '  -3.88%  ' },
    @{ Col = 'D'; Row = 15; Value = '7.045' },
    @{ Col = 'E'; Row = 15; Value = '  -2.92%  ' },
    @{ Col = 'D'; Row = 16; Value = '1.006' },
    @{ Col = 'E'; Row = 16; Value = '  -1.08%  ' },
    @{ Col = 'D'; Row = 17; Value = '87.98' },
    @{ Col = 'E'; Row = 17; Value = '  -1.10%  ' },
    @{ Col = 'D'; Row = 18; Value = '0.06635' },
    @{ Col = 'E'; Row = 18; Value = '  -2.14%  ' },
    @{ Col = 'D'; Row = 19; Value = '0.00001031' },
    @{ Col = 'E'; Row = 19; Value = '  -2.33%  ' },
    @{ Col = 'D'; Row = 20; Value = '16.98' },
    @{ Col = 'E'; Row = 20; Value = '  -1.19%  ' },
    @{ Col = 'D'; Row = 21; Value = '1.005' },
    @{ Col = 'E'; Row = 21; Value = '  -0.90%  ' },
    @{ Col = 'D'; Row = 22; Value = '27.512.17' },
    @{ Col = 'E'; Row = 22; Value = '  -2.10%  ' },
    @{ Col = 'D'; Row = 23; Value = '5.355' },
    @{ Col = 'E'; Row = 23; Value = '  -3.68%  ' },
    @{ Col = 'D'; Row = 24; Value = '10.90' },
    @{ Col = 'E'; Row = 24; Value = '  -1.73%  ' },
    @{ Col = 'D'; Row = 25; Value = '2.295' },
    @{ Col = 'E'; Row = 25; Value = '  -2.60%  ' },
    @{ Col = 'D'; Row = 26; Value = '157.13' },
    @{ Col = 'E'; Row = 26; Value = '  -2.45%  ' },
    @{ Col = 'D'; Row = 27; Value = '19.41' },
    @{ Col = 'E'; Row = 27; Value = '  -3.73%  ' },
    @{ Col = 'D'; Row = 28; Value = '2.078' },
    @{ Col = 'E'; Row = 28; Value = '  -2.26%  ' },
    @{ Col = 'D'; Row = 29; Value = '5.356' },
    @{ Col = 'E'; Row = 29; Value = '  -3.93%  ' },
    @{ Col = 'D'; Row = 30; Value = '118.91' },
    @{ Col = 'E'; Row = 30; Value = '  -2.72%  ' },
    @{ Col = 'D'; Row = 31; Value = '0.9585' },
    @{ Col = 'E'; Row = 31; Value = '  -2.68%  ' },
    @{ Col = 'D'; Row = 32; Value = '0.09331' },
    @{ Col = 'E'; Row = 32; Value = '  -3.01%  ' },
    @{ Col = 'D'; Row = 33; Value = '3.573' },
    @{ Col = 'E'; Row = 33; Value = '  -2.13%  ' },
    @{ Col = 'D'; Row = 34; Value = '5.246' },
    @{ Col = 'E'; Row = 34; Value = '  -2.64%  ' },
    @{ Col = 'D'; Row = 35; Value = '1.325' },
    @{ Col = 'E'; Row = 35; Value = '  -3.96%  ' },
    @{ Col = 'B'; Row = 36; Value = 'Hedera' },
    @{ Col = 'C'; Row = 36; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' },
    @{ Col = 'D'; Row = 36; Value = '0.05949' },
    @{ Col = 'E'; Row = 36; Value = '  -2.95%  ' },
    @{ Col = 'B'; Row = 37; Value = 'VeChain' },
    @{ Col = 'C'; Row = 37; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' },
    @{ Col = 'D'; Row = 37; Value = '0.02204' },
    @{ Col = 'E'; Row = 37; Value = '  -2.76%  ' },
    @{ Col = 'D'; Row = 38; Value = '8.102' },
    @{ Col = 'E'; Row = 38; Value = '  -1.80%  ' },
    @{ Col = 'D'; Row = 39; Value = '1.161' },
    @{ Col = 'E'; Row = 39; Value = '  -3.82%  ' },
    @{ Col = 'D'; Row = 40; Value = '0.5847' },
    @{ Col = 'E'; Row = 40; Value = '  -2.67%  ' },
    @{ Col = 'D'; Row = 41; Value = '0.1847' },
    @{ Col = 'E'; Row = 41; Value = '  -3.08%  ' },
    @{ Col = 'E'; Row = 42; Value = '  -2.78%  ' },
    @{ Col = 'D'; Row = 43; Value = '1.258' },
    @{ Col = 'E'; Row = 43; Value = '  -1.81%  ' },
    @{ Col = 'D'; Row = 44; Value = '0.5523' },
    @{ Col = 'E'; Row = 44; Value = '  -3.13%  ' },
    @{ Col = 'D'; Row = 45; Value = '12.03' },
    @{ Col = 'E'; Row = 45; Value = '  -2.10%  ' },
    @{ Col = 'D'; Row = 46; Value = '1.871' },
    @{ Col = 'E'; Row = 46; Value = '  -3.79%  ' },
    @{ Col = 'D'; Row = 47; Value = '0.06660' },
    @{ Col = 'E'; Row = 47; Value = '  -2.70%  ' },
    @{ Col = 'D'; Row = 48; Value = '110.48' },
    @{ Col = 'E'; Row = 48; Value = '  -2.07%  ' },
    @{ Col = 'D'; Row = 49; Value = '1.045' },
    @{ Col = 'E'; Row = 49; Value = '  -3.19%  ' },
    @{ Col = 'E'; Row = 50; Value = '  -2.38%  ' },
    @{ Col = 'D'; Row = 51; Value = '1.005' },
    @{ Col = 'E'; Row = 51; Value = '  -1.22%  ' }
)

foreach ($change in $changes) {
    $cellRef = "$($change.Col)$($change.Row)"
    $cell = $ws.Range($cellRef)

    if ($change.Col -eq 'D') {
        # Many "Price" values look numeric (e.g. "1.005", "0.00001031"); force them to be
        # stored as plain text (matching the source data) via a leading apostrophe, then
        # strip the resulting "quote prefix" formatting so the cell style is left untouched.
        $cell.Value = "'" + $change.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $change.Value
    }
}

Write-Host "Applied $($changes.Count) cell updates"
